# The commit regenerates the document's list-numbering definitions, which
# causes Word to mint fresh `nsid` (list "name space" / GUID-like) values
# for the abstractNum entries it rewrote. The `w:nsid` attribute is an
# internal OOXML identifier that isn't exposed anywhere on the Word object
# model (ListTemplates/ListLevels, etc. don't surface it), so we go
# straight at the package's flat-OPC XML via Document.WordOpenXML, patch
# the four nsid values for the touched abstractNum definitions, and write
# the XML back.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('<w:nsid w:val="41181277" />', '<w:nsid w:val="480a704c" />')
$xml = $xml.Replace('<w:nsid w:val="cdcf4788" />', '<w:nsid w:val="803dfc18" />')
$xml = $xml.Replace('<w:nsid w:val="6fdbc4d1" />', '<w:nsid w:val="b8f76748" />')
$xml = $xml.Replace('<w:nsid w:val="9275a25c" />', '<w:nsid w:val="ecc12655" />')

$d.WordOpenXML = $xml
